$d = $word.ActiveDocument

# Update the date line (first paragraph, outside the table).
$d.Paragraphs.Item(1).Range.Text = "2023-08-27 Sunday"

# Update the table of division problems. Each data row/col is addressed
# directly via Table.Cell(row, col) and the cell's Range.Text is assigned
# in place (after trimming the trailing end-of-cell mark). This keeps each
# edit confined to exactly the intended cell - several cells in this table
# share identical text, so a document-wide Find/Replace would risk touching
# the wrong cell.
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cellRange = $table.Cell($row, $col).Range
    $cellRange.MoveEnd(1, -1) | Out-Null
    $cellRange.Text = $newText
}

Set-CellText $t 1 1 "95÷8=11, 7"
Set-CellText $t 1 2 "62÷2=31, 0"
Set-CellText $t 1 3 "52÷8=6, 4"
Set-CellText $t 1 4 "69÷5=13, 4"
# Row 1 Col 5 "63÷6=10, 3" is unchanged

Set-CellText $t 5 1 "91÷2=45, 1"
Set-CellText $t 5 2 "99÷9=11, 0"
Set-CellText $t 5 3 "79÷7=11, 2"
Set-CellText $t 5 4 "55÷2=27, 1"
Set-CellText $t 5 5 "87÷4=21, 3"

Set-CellText $t 9 1 "74÷4=18, 2"
Set-CellText $t 9 2 "66÷8=8, 2"
Set-CellText $t 9 3 "79÷2=39, 1"
Set-CellText $t 9 4 "12÷5=2, 2"
Set-CellText $t 9 5 "28÷6=4, 4"

Set-CellText $t 13 1 "66÷8=8, 2"
Set-CellText $t 13 2 "35÷7=5, 0"
Set-CellText $t 13 3 "84÷7=12, 0"
Set-CellText $t 13 4 "56÷4=14, 0"
Set-CellText $t 13 5 "75÷7=10, 5"

Set-CellText $t 17 1 "12÷8=1, 4"
Set-CellText $t 17 2 "35÷2=17, 1"
Set-CellText $t 17 3 "26÷7=3, 5"
Set-CellText $t 17 4 "44÷7=6, 2"
Set-CellText $t 17 5 "20÷2=10, 0"
